$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 383.33334
$ws.Range("I33").Value = 300.0
$ws.Range("J33").Value = 800.0
$ws.Range("K33").Value = 300.0
$ws.Range("L33").Value = 800.0
$ws.Range("M33").Value = -71.0
$ws.Range("N33").Value = -1258.0
$ws.Range("H111").Value = 3000.0
$ws.Range("I111").Value = 0.0
$ws.Range("J111").Value = 3000.0
$ws.Range("K111").Value = 0.0
$ws.Range("L111").Value = 9000.0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -15134.0
$ws.Range("H129").Value = 1366.9854
$ws.Range("J129").Value = 1454.4354
$ws.Range("L129").Value = 4363.3062
$ws.Range("N129").Value = -14363.3062
$ws.Range("H132").Value = 26258736.0
$ws.Range("I132").Value = 30729244.0
$ws.Range("J132").Value = 1670944.1
$ws.Range("K132").Value = 92187732.0
$ws.Range("L132").Value = 5012832.300000001
$ws.Range("M132").Value = -92185202.0
$ws.Range("N132").Value = -5017892.300000001
$ws.Range("H137").Value = 838866.75
$ws.Range("I137").Value = 2168500.5
$ws.Range("J137").Value = 3097.0286
$ws.Range("K137").Value = 6505501.5
$ws.Range("L137").Value = 9291.0858
$ws.Range("M137").Value = -6502951.5
$ws.Range("N137").Value = -14391.0858
$ws.Range("H138").Value = 2995.3845
$ws.Range("J138").Value = 3918.75
$ws.Range("L138").Value = 11756.25
$ws.Range("N138").Value = -22036.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3833.0322
$ws.Range("I122").Value = 3628.4092
$ws.Range("J122").Value = 4333.222
$ws.Range("K122").Value = 10885.2276
$ws.Range("L122").Value = 12999.666
$ws.Range("M122").Value = -8435.2276
$ws.Range("N122").Value = -17899.666
$ws.Range("H132").Value = 2835.5925
$ws.Range("I132").Value = 2302.3333
$ws.Range("J132").Value = 4702.0
$ws.Range("K132").Value = 6906.999899999999
$ws.Range("L132").Value = 14106.0
$ws.Range("M132").Value = -4376.999899999999
$ws.Range("N132").Value = -19166.0
$ws.Range("H137").Value = 39538.5
$ws.Range("J137").Value = 39538.5
$ws.Range("L137").Value = 39538.5
$ws.Range("N137").Value = -49738.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2926.52
$ws.Range("I134").Value = 1038.909
$ws.Range("J134").Value = 6590.706
$ws.Range("K134").Value = 3116.727
$ws.Range("L134").Value = 19772.118
$ws.Range("M134").Value = -581.7270000000003
$ws.Range("N134").Value = -24842.118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 243799.02
$ws.Range("I31").Value = 966354.4
$ws.Range("J31").Value = 2947.238
$ws.Range("K31").Value = 966354.4
$ws.Range("L31").Value = 2947.238
$ws.Range("M31").Value = -966059.4
$ws.Range("N31").Value = -3537.238
$ws.Range("H34").Value = 243799.02
$ws.Range("I34").Value = 966354.4
$ws.Range("J34").Value = 2947.238
$ws.Range("K34").Value = 966354.4
$ws.Range("L34").Value = 2947.238
$ws.Range("M34").Value = -966152.4
$ws.Range("N34").Value = -3351.238
$ws.Range("H98").Value = 47800.0
$ws.Range("J98").Value = 47800.0
$ws.Range("L98").Value = 47800.0
$ws.Range("N98").Value = -52292.0
$ws.Range("H141").Value = 27232.0
$ws.Range("J141").Value = 27232.0
$ws.Range("L141").Value = 27232.0
$ws.Range("N141").Value = -37592.0
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 557172.5
$ws.Range("I5").Value = 509.6
$ws.Range("J5").Value = 703662.75
$ws.Range("K5").Value = 1528.8
$ws.Range("L5").Value = 2110988.25
$ws.Range("M5").Value = -1416.8
$ws.Range("N5").Value = -2111212.25
$ws.Range("H68").Value = 3812.907
$ws.Range("I68").Value = 1093.875
$ws.Range("J68").Value = 5424.185
$ws.Range("K68").Value = 3281.625
$ws.Range("L68").Value = 16272.555
$ws.Range("M68").Value = -2470.625
$ws.Range("N68").Value = -17894.555
$ws.Range("H71").Value = 3812.907
$ws.Range("I71").Value = 1093.875
$ws.Range("J71").Value = 5424.185
$ws.Range("K71").Value = 9844.875
$ws.Range("L71").Value = 48817.665
$ws.Range("M71").Value = -5788.875
$ws.Range("N71").Value = -56929.665
$ws.Range("H92").Value = 799.25
$ws.Range("I92").Value = 798.0
$ws.Range("J92").Value = 803.0
$ws.Range("K92").Value = 2394.0
$ws.Range("L92").Value = 2409.0
$ws.Range("M92").Value = -1146.0
$ws.Range("N92").Value = -4905.0
$ws.Range("H113").Value = 1689876.9
$ws.Range("I113").Value = 631.9808
$ws.Range("J113").Value = 5682637.5
$ws.Range("K113").Value = 1895.9424
$ws.Range("L113").Value = 17047912.5
$ws.Range("M113").Value = 274.0575999999999
$ws.Range("N113").Value = -17052252.5
$ws.Range("H131").Value = 862.71
$ws.Range("I131").Value = 538.6
$ws.Range("J131").Value = 879.76843
$ws.Range("K131").Value = 1615.8
$ws.Range("L131").Value = 2639.30529
$ws.Range("M131").Value = 3424.2
$ws.Range("N131").Value = -12719.30529
$ws.Range("H135").Value = 557172.5
$ws.Range("I135").Value = 509.6
$ws.Range("J135").Value = 703662.75
$ws.Range("K135").Value = 4586.400000000001
$ws.Range("L135").Value = 6332964.75
$ws.Range("M135").Value = -2051.400000000001
$ws.Range("N135").Value = -6338034.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28069.23
$ws.Range("J15").Value = 28069.23
$ws.Range("L15").Value = 28069.23
$ws.Range("N15").Value = -28645.23
$ws.Range("H81").Value = 28069.23
$ws.Range("J81").Value = 28069.23
$ws.Range("L81").Value = 28069.23
$ws.Range("N81").Value = -30065.23
$ws.Range("H84").Value = 28069.23
$ws.Range("J84").Value = 28069.23
$ws.Range("L84").Value = 84207.69
$ws.Range("N84").Value = -94191.69
$ws.Range("H102").Value = 6199.8
$ws.Range("I102").Value = 5083.3335
$ws.Range("J102").Value = 7874.5
$ws.Range("K102").Value = 5083.3335
$ws.Range("L102").Value = 7874.5
$ws.Range("M102").Value = -3461.3335
$ws.Range("N102").Value = -11118.5
$ws.Range("H126").Value = 4052.0112
$ws.Range("I126").Value = 2998.255
$ws.Range("J126").Value = 5430.0
$ws.Range("K126").Value = 8994.765
$ws.Range("L126").Value = 16290.0
$ws.Range("M126").Value = -6524.764999999999
$ws.Range("N126").Value = -21230.0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6252.143
$ws.Range("I40").Value = 5684.5454
$ws.Range("K40").Value = 5684.5454
$ws.Range("M40").Value = -5548.5454
$ws.Range("H46").Value = 1033.625
$ws.Range("I46").Value = 680.7273
$ws.Range("K46").Value = 680.7273
$ws.Range("M46").Value = -492.7273
$ws.Range("H100").Value = 2980.8
$ws.Range("I100").Value = 2100.0
$ws.Range("J100").Value = 3201.0
$ws.Range("K100").Value = 2100.0
$ws.Range("L100").Value = 3201.0
$ws.Range("M100").Value = -1559.0
$ws.Range("N100").Value = -4283.0
$ws.Range("H118").Value = 0.0
$ws.Range("J118").Value = 0.0
$ws.Range("L118").Value = 0.0
$ws.Range("N118").ClearContents()
$ws.Range("H123").Value = 43827.0
$ws.Range("J123").Value = 43827.0
$ws.Range("L123").Value = 43827.0
$ws.Range("N123").Value = -53627.0
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 72780.625
$ws.Range("I62").Value = 3668.4
$ws.Range("J62").Value = 187967.67
$ws.Range("K62").Value = 3668.4
$ws.Range("L62").Value = 187967.67
$ws.Range("M62").Value = -3044.4
$ws.Range("N62").Value = -189215.67
$ws.Range("H65").Value = 72780.625
$ws.Range("I65").Value = 3668.4
$ws.Range("J65").Value = 187967.67
$ws.Range("K65").Value = 18342.0
$ws.Range("L65").Value = 939838.3500000001
$ws.Range("M65").Value = -15222.0
$ws.Range("N65").Value = -946078.3500000001
$ws.Range("H107").Value = 855.0
$ws.Range("I107").Value = 609.0
$ws.Range("K107").Value = 1827.0
$ws.Range("M107").Value = 93.0
$ws.Range("H122").Value = 3629.6553
$ws.Range("I122").Value = 2204.7273
$ws.Range("K122").Value = 6614.1819
$ws.Range("M122").Value = -4164.1819
